$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-18 23:18:45"
$ws.Range("H2").Value = "'70%"
$ws.Range("I2").Value = "1.5 mm"
$ws.Range("E3").Value = "2026-02-18 23:18:48"
$ws.Range("O3").Value = "-0.3 °C"
$ws.Range("E4").Value = "2026-02-18 23:18:51"
$ws.Range("J4").Value = "1012.1 hPa"
$ws.Range("K4").Value = "11.7 MJ/m2"
$ws.Range("E5").Value = "2026-02-18 23:18:54"
$ws.Range("G5").Value = "137 cm"
$ws.Range("H5").Value = "'69%"
$ws.Range("I5").Value = "1.8 mm"
$ws.Range("N5").Value = "-4.6 °C 22:59 TU"
$ws.Range("O5").Value = "0.4 °C"
$ws.Range("E6").Value = "2026-02-18 23:18:56"
$ws.Range("J6").Value = "1011.8 hPa"
$ws.Range("E7").Value = "2026-02-18 23:18:59"
$ws.Range("J7").Value = "1013.4 hPa"
$ws.Range("E8").Value = "2026-02-18 23:19:02"
$ws.Range("H8").Value = "'82%"
$ws.Range("J8").Value = "1013.2 hPa"
$ws.Range("L8").Value = "47.5 km/h - 230º 22:58 TU"
$ws.Range("E9").Value = "2026-02-18 23:19:04"
$ws.Range("H9").Value = "'76%"
$ws.Range("O9").Value = "11.2 °C"
$ws.Range("E10").Value = "2026-02-18 23:19:06"
$ws.Range("H10").Value = "'84%"
$ws.Range("O10").Value = "10.7 °C"
$ws.Range("E11").Value = "2026-02-18 23:19:09"
$ws.Range("E12").Value = "2026-02-18 23:19:12"
$ws.Range("H12").Value = "'84%"
$ws.Range("E13").Value = "2026-02-18 23:19:14"
$ws.Range("H13").Value = "'74%"
$ws.Range("J13").Value = "1014.4 hPa"
$ws.Range("E14").Value = "2026-02-18 23:19:17"
$ws.Range("O14").Value = "12.2 °C"
$ws.Range("E15").Value = "2026-02-18 23:19:20"
$ws.Range("E16").Value = "2026-02-18 23:19:22"
$ws.Range("N16").Value = "-4.4 °C 22:59 TU"
$ws.Range("O16").Value = "-0.4 °C"
$ws.Range("E17").Value = "2026-02-18 23:19:25"
$ws.Range("O17").Value = "3.4 °C"
$ws.Range("E18").Value = "2026-02-18 23:19:28"
$ws.Range("J18").Value = "1012.3 hPa"
$ws.Range("O18").Value = "11.9 °C"
$ws.Range("E19").Value = "2026-02-18 23:19:31"
$ws.Range("E20").Value = "2026-02-18 23:19:33"
$ws.Range("L20").Value = "57.6 km/h - 285º 22:32 TU"
$ws.Range("N20").Value = "-4.2 °C 22:59 TU"
$ws.Range("E21").Value = "2026-02-18 23:19:36"
$ws.Range("J21").Value = "1013.9 hPa"
$ws.Range("E22").Value = "2026-02-18 23:19:39"
$ws.Range("E23").Value = "2026-02-18 23:19:42"
$ws.Range("N23").Value = "-3.9 °C 22:57 TU"
$ws.Range("E24").Value = "2026-02-18 23:19:44"
$ws.Range("J24").Value = "1014.1 hPa"
$ws.Range("E25").Value = "2026-02-18 23:19:47"
$ws.Range("I25").Value = "0.2 mm"
$ws.Range("E26").Value = "2026-02-18 23:19:49"
$ws.Range("J26").Value = "1011.4 hPa"
$ws.Range("E27").Value = "2026-02-18 23:19:52"
$ws.Range("H27").Value = "'60%"
$ws.Range("N27").Value = "-2.3 °C 22:59 TU"
$ws.Range("E28").Value = "2026-02-18 23:19:55"
$ws.Range("J28").Value = "1012.0 hPa"
$ws.Range("E29").Value = "2026-02-18 23:19:58"
$ws.Range("E30").Value = "2026-02-18 23:20:01"
$ws.Range("J30").Value = "1011.6 hPa"
$ws.Range("E31").Value = "2026-02-18 23:20:03"
$ws.Range("J31").Value = "1010.5 hPa"
$ws.Range("E32").Value = "2026-02-18 23:20:06"
$ws.Range("E33").Value = "2026-02-18 23:20:09"
$ws.Range("H33").Value = "'68%"
$ws.Range("J33").Value = "1013.2 hPa"
$ws.Range("E34").Value = "2026-02-18 23:20:12"
$ws.Range("H34").Value = "'50%"
$ws.Range("O34").Value = "2.9 °C"
$ws.Range("E35").Value = "2026-02-18 23:20:15"
$ws.Range("I35").Value = "0.5 mm"
$ws.Range("O35").Value = "9.1 °C"
$ws.Range("E36").Value = "2026-02-18 23:20:17"
$ws.Range("H36").Value = "'83%"
$ws.Range("J36").Value = "1012.1 hPa"
$ws.Range("E37").Value = "2026-02-18 23:20:20"
$ws.Range("J37").Value = "1013.7 hPa"
$ws.Range("E38").Value = "2026-02-18 23:20:23"
$ws.Range("K38").Value = "12.7 MJ/m2"
$ws.Range("E39").Value = "2026-02-18 23:20:25"
$ws.Range("H39").Value = "'44%"
$ws.Range("I39").Value = "0.3 mm"
$ws.Range("O39").Value = "1.0 °C"
$ws.Range("E40").Value = "2026-02-18 23:20:28"
$ws.Range("H40").Value = "'78%"
$ws.Range("J40").Value = "1014.6 hPa"
$ws.Range("E41").Value = "2026-02-18 23:20:31"
$ws.Range("J41").Value = "1013.8 hPa"
$ws.Range("L41").Value = "34.9 km/h - 247º 22:57 TU"
$ws.Range("E42").Value = "2026-02-18 23:20:33"
$ws.Range("O42").Value = "11.9 °C"
$ws.Range("E43").Value = "2026-02-18 23:20:36"
$ws.Range("E44").Value = "2026-02-18 23:20:39"
$ws.Range("E45").Value = "2026-02-18 23:20:41"
$ws.Range("O45").Value = "7.1 °C"
$ws.Range("E46").Value = "2026-02-18 23:20:44"
$ws.Range("J46").Value = "1014.2 hPa"
$ws.Range("L46").Value = "31.0 km/h - 308º 22:47 TU"
